$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells whose new values look like plain numbers stay as text,
# matching the workbook convention of storing Price/Volume as strings.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.213.03"
$ws.Range("E2").Value = "  +0.93%  "
$ws.Range("D3").Value = "1.570.16"
$ws.Range("E3").Value = "  +0.61%  "
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("D5").Value = "211.03"
$ws.Range("E5").Value = "  +1.85%  "
$ws.Range("E6").Value = "  +0.45%  "
$ws.Range("E7").Value = "  +0.29%  "
$ws.Range("D8").Value = "22.05"
$ws.Range("E8").Value = "  -0.36%  "
$ws.Range("E9").Value = "  +0.34%  "
$ws.Range("E10").Value = "  +0.08%  "
$ws.Range("D11").Value = "0.0867"
$ws.Range("E11").Value = "  +1.10%  "
$ws.Range("D12").Value = "1.793.93"
$ws.Range("E12").Value = "  +0.60%  "
$ws.Range("D13").Value = "1.599.92"
$ws.Range("E13").Value = "  +2.35%  "
$ws.Range("E14").Value = "  +0.90%  "
$ws.Range("E15").Value = "  -0.02%  "
$ws.Range("D16").Value = "27.163.49"
$ws.Range("E16").Value = "  +0.70%  "
$ws.Range("D17").Value = "62.26"
$ws.Range("D18").Value = "7.51"
$ws.Range("E18").Value = "  +2.01%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "216.53"
$ws.Range("E19").Value = "  -0.32%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0703"
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("E21").Value = "  +0.32%  "
$ws.Range("E22").Value = "  +1.02%  "
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("E24").Value = "  +0.30%  "
$ws.Range("D25").Value = "154.10"
$ws.Range("E25").Value = "  +0.46%  "
$ws.Range("D26").Value = "6.64"
$ws.Range("E26").Value = "  +0.40%  "
$ws.Range("D27").Value = "15.07"
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").Value = "0.107"
$ws.Range("E28").Value = "  +2.04%  "
$ws.Range("E29").Value = "  +0.25%  "
$ws.Range("E30").Value = "  +2.71%  "
$ws.Range("E31").Value = "  +0.64%  "
$ws.Range("D32").Value = "3.24"
$ws.Range("E32").Value = "  +0.16%  "
$ws.Range("D33").Value = "1.450.19"
$ws.Range("E33").Value = "  +1.98%  "
$ws.Range("E34").Value = "  +1.35%  "
$ws.Range("E35").Value = "  +6.93%  "
$ws.Range("E36").Value = "  +0.51%  "
$ws.Range("D37").Value = "2.34"
$ws.Range("E37").Value = "  +0.59%  "
$ws.Range("E38").Value = "  +1.08%  "
$ws.Range("D39").Value = "0.533"
$ws.Range("E39").Value = "  -0.33%  "
$ws.Range("D40").Value = "5.85"
$ws.Range("E40").Value = "  +2.29%  "
$ws.Range("E41").Value = "  +0.23%  "
$ws.Range("E42").Value = "  +0.34%  "
$ws.Range("E43").Value = "  +0.37%  "
$ws.Range("E44").Value = "  -0.24%  "
$ws.Range("E45").Value = "  -0.67%  "
$ws.Range("D46").Value = "1.72"
$ws.Range("E46").Value = "  -1.23%  "
$ws.Range("D47").Value = "1.705.38"
$ws.Range("E47").Value = "  +0.45%  "
$ws.Range("D48").Value = "86.11"
$ws.Range("E48").Value = "  -1.50%  "
$ws.Range("E49").Value = "  +3.64%  "
$ws.Range("E50").Value = "  +0.40%  "
$ws.Range("D51").Value = "0.0957"
$ws.Range("E51").Value = "  +0.03%  "
